$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.741.89"
$ws.Range("E2").Value = "  +4.02%  "

$ws.Range("D3").Value = "2.665.05"
$ws.Range("E3").Value = "  +7.52%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "325.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.99%  "

$ws.Range("E11").Value = "  -1.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0822"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.73%  "

$ws.Range("D15").Value = "3.090.30"
$ws.Range("E15").Value = "  +7.74%  "

$ws.Range("D16").Value = "2.680.47"
$ws.Range("E16").Value = "  +7.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.873"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.32%  "

$ws.Range("D18").Value = "49.715.20"
$ws.Range("E18").Value = "  +4.12%  "

$ws.Range("E19").Value = "  +4.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.18%  "

$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("D22").Value = "0.0₃0956"
$ws.Range("E22").Value = "  +3.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "279.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.87%  "

$ws.Range("E25").Value = "  +3.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.87%  "

$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.84%  "

$ws.Range("E31").Value = "  +4.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.99%  "

$ws.Range("E33").Value = "  +4.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0811"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.69%  "

$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.61%  "

$ws.Range("E41").Value = "  +2.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.16%  "

$ws.Range("E43").Value = "  +1.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0317"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.61%  "

$ws.Range("D45").Value = "2.106.11"
$ws.Range("E45").Value = "  +5.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.75%  "

$ws.Range("E47").Value = "  +14.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.38%  "
